$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.499.87"
$ws.Range("E2").Value = "  +5.39%  "
$ws.Range("D3").Value = "1.724.53"
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'225.32"
$ws.Range("E5").Value = "  +3.02%  "
$ws.Range("E6").Value = "  +2.82%  "
$ws.Range("D7").Value = "'1.003"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.2661"
$ws.Range("E8").Value = "  +1.47%  "
$ws.Range("D9").Value = "'0.06605"
$ws.Range("E9").Value = "  +4.81%  "
$ws.Range("D10").Value = "'21.56"
$ws.Range("E10").Value = "  +6.15%  "
$ws.Range("D11").Value = "'0.07664"
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").Value = "'4.599"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").Value = "1.724.51"
$ws.Range("E13").Value = "  +3.61%  "
$ws.Range("D14").Value = "1.961.77"
$ws.Range("E14").Value = "  +4.63%  "
$ws.Range("D15").Value = "'0.5803"
$ws.Range("E15").Value = "  +4.13%  "
$ws.Range("D16").Value = "0.0₅8293"
$ws.Range("E16").Value = "  +2.18%  "
$ws.Range("D17").Value = "'67.82"
$ws.Range("E17").Value = "  +4.22%  "
$ws.Range("D18").Value = "27.493.65"
$ws.Range("E18").Value = "  +5.49%  "
$ws.Range("D19").Value = "'218.63"
$ws.Range("E19").Value = "  +12.48%  "
$ws.Range("D20").Value = "'1.003"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("E21").Value = "  +2.77%  "
$ws.Range("D22").Value = "'10.59"
$ws.Range("E22").Value = "  +1.31%  "
$ws.Range("D23").Value = "'6.035"
$ws.Range("E23").Value = "  +1.96%  "
$ws.Range("D24").Value = "'1.003"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").Value = "'143.06"
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("D26").Value = "'1.753"
$ws.Range("E26").Value = "  +15.95%  "
$ws.Range("D27").Value = "'0.1233"
$ws.Range("E27").Value = "  +4.40%  "
$ws.Range("D28").Value = "'7.328"
$ws.Range("E28").Value = "  +1.90%  "
$ws.Range("E29").Value = "  +3.99%  "
$ws.Range("D30").Value = "'0.05490"
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("D31").Value = "'1.300"
$ws.Range("E31").Value = "  +2.46%  "
$ws.Range("E32").Value = "  +3.49%  "
$ws.Range("D33").Value = "'3.442"
$ws.Range("E33").Value = "  +3.54%  "
$ws.Range("D34").Value = "'1.661"
$ws.Range("E34").Value = "  +6.74%  "
$ws.Range("D35").Value = "'2.856"
$ws.Range("E35").Value = "  +2.72%  "
$ws.Range("D36").Value = "'0.9571"
$ws.Range("E36").Value = "  +1.67%  "
$ws.Range("D37").Value = "'2.423"
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("D38").Value = "'0.5936"
$ws.Range("E38").Value = "  +5.95%  "
$ws.Range("D39").Value = "'0.01649"
$ws.Range("E39").Value = "  +4.99%  "
$ws.Range("D40").Value = "'5.899"
$ws.Range("E40").Value = "  +2.89%  "
$ws.Range("D41").Value = "1.049.35"
$ws.Range("E41").Value = "  +2.09%  "
$ws.Range("D42").Value = "'0.8471"
$ws.Range("E42").Value = "  +2.88%  "
$ws.Range("D43").Value = "'1.003"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "'101.36"
$ws.Range("E44").Value = "  +0.71%  "
$ws.Range("D45").Value = "1.868.01"
$ws.Range("E45").Value = "  +4.59%  "
$ws.Range("E46").Value = "  +4.54%  "
$ws.Range("D47").Value = "'58.79"
$ws.Range("E47").Value = "  +2.64%  "
$ws.Range("E48").Value = "  +3.66%  "
$ws.Range("D49").Value = "'8.215"
$ws.Range("E49").Value = "  +4.31%  "
$ws.Range("D50").Value = "'1.003"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").Value = "'0.05247"
$ws.Range("E51").Value = "  +2.90%  "
